$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.156.06'
$ws.Range('E2').Value = '  +0.15%  '
$ws.Range('D3').Value = '2.913.99'
$ws.Range('E3').Value = '  +0.37%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '369.55'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +6.46%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '103.71'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.45%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.539'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.49%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.585'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.74%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.59'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.84%  '
$ws.Range('E11').Value = '  +1.23%  '
$ws.Range('E12').Value = '  -1.05%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.33'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.15%  '
$ws.Range('D14').Value = '3.370.37'
$ws.Range('E14').Value = '  +0.23%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.39'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.66%  '
$ws.Range('D16').Value = '2.914.40'
$ws.Range('E16').Value = '  -2.01%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.933'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.02%  '
$ws.Range('D18').Value = '51.138.37'
$ws.Range('E18').Value = '  +0.29%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.23'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.15%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.19'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.91%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.92'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.63%  '
$ws.Range('D22').Value = '0.0₃0943'
$ws.Range('E22').Value = '  -1.26%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '68.42'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.06%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '259.45'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.07%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.69'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.00%  '
$ws.Range('E26').Value = '  +3.92%  '
$ws.Range('E27').Value = '  +4.28%  '
$ws.Range('E28').Value = '  -0.10%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '25.71'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.61%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.04'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.83%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.102'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.16%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.23'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.35%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '9.88'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.31%  '
$ws.Range('E34').Value = '  +0.77%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '34.57'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.31%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '50.86'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.43%  '
$ws.Range('E37').Value = '  +0.56%  '
$ws.Range('E38').Value = '  +0.68%  '
$ws.Range('E39').Value = '  -2.71%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.65'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.77%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '17.08'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.05%  '
$ws.Range('E42').Value = '  -3.99%  '
$ws.Range('E43').Value = '  -1.86%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '22.29'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.19%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '118.34'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.68%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.08'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.08%  '
$ws.Range('D47').Value = '2.016.36'
$ws.Range('E47').Value = '  -3.08%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.16'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.25%  '
$ws.Range('D50').Value = '3.210.82'
$ws.Range('E50').Value = '  +0.62%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.240'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.93%  '
